$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    1086.004175901413,
    214.3675310611725,
    909.3178842067719,
    675.5533845424652,
    279.4213538169861,
    606.6833273172379,
    395.1025931835175,
    1385.716940760612,
    1327.859076142311,
    1917.706836104393,
    941.903670668602,
    123.012192606926,
    678.4823168516159,
    53.69848430156708,
    67.35533845424652,
    75.25483322143555,
    226.1076455116272,
    296.3919162750244,
    177.622364282608,
    191.3380935192108,
    1566.243852138519,
    455.9726504087448,
    226.0071402788162,
    374.2741661071777,
    677.5533845424652,
    579.3107439279556,
    581.4112491607666,
    965.1412589550018,
    163.2375882863998,
    244.8771975040436,
    479.9482651948929,
    322.7766922712326,
    284.634556889534,
    319.1198381185532,
    198.4507913589478,
    826.7493426799774,
    203.8650048971176,
    114.5685415267944,
    237.0365778207779,
    220.7939372062683,
    470.9432128667831,
    96.56854152679443,
    60.87005722522736,
    495.9137753248215,
    109.6396092176437,
    276.1492756605148,
    217.0365778207779,
    58.18376553058624,
    355.64674949646,
    143.923879981041,
    475.2863587141037,
    600.180801153183,
    820.5483322143555,
    71.84061968326569,
    125.2964633703232,
    383.6883796453476,
    243.7056245803833,
    296.0487704277039,
    297.3624787330627,
    1068.864128470421,
    556.7422024011612,
    482.0142805576324,
    237.9066350460052,
    659.5950146913528,
    329.2619735002518,
    334.1320307254791,
    281.2203433513641,
    429.7888848781586,
    212.3502861261368,
    458.6000670194626,
    508.8843377828598,
    253.8061298131943,
    156.2670258283615,
    293.9482651948929,
    219.0365778207779,
    201.9655101299286,
    182.4091612100601,
    260.1492756605148,
    532.3401814699173,
    865.5432798862457,
    1108.454967260361,
    240.4924215078354,
    503.7716399431229,
    598.0975408554077,
    381.8893901109695,
    703.1513636112213,
    209.82337474823,
    553.9137753248215,
    693.4945094585419,
    300.2914110422134,
    248.5929267406464,
    207.8650048971176,
    286.3919162750244,
    693.2518688440323,
    145.5807341337204,
    299.0193328857422,
    866.0874361991882,
    382.2325359582901,
    586.2396762371063,
    311.9482651948929,
    456.7594473361969,
    384.4507913589478,
    194.4507913589478,
    499.0437180995941,
    174.7523070573807,
    530.3990565538406,
    650.9675980806351,
    281.5634891986847,
    413.8721451759338,
    56.87005722522736,
    220.1665205955505,
    243.7644996643066,
    171.82337474823,
    263.4213538169861,
    203.4802289009094,
    123.9827550649643,
    74.04163014888763,
    160.0243852138519
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $values[$i]
}
